$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(8.044522085013483, 4.051576211009759, -9.405187126396871, 10.50532368479741)"
$ws.Range("C2").Value = "NIG(0.9681324837651086, 0.7278351810495378, 4.6871252690949605, 5.627473558382173)"
$ws.Range("D2").Value = "JSU(-0.8284919241294476, 1.0230501081776773, 0.997942388522383, 2.391106250688419)"
$ws.Range("E2").Value = "NCT(3.5004362600298995, 1.91752511366082, -0.018546684878672854, 4.347317552970066)"
